# Update the "想去人数" (want-to-go count) figures in column F for both the
# "展览" and "全部类型" worksheets (they hold duplicated data).
#
# Mapping of row -> old/new value (from F2..F21):
#   F2 : 11504 -> 11526
#   F3 : 10981 -> 11008
#   F6 :   994 -> 998
#   F7 :   116 -> 117
#   F9 :    39 -> 40
#   F10:    37 -> 38
#   F11: 10628 -> 10640
#   F12:  4109 -> 4113
#   F15:  2456 -> 2457
#   F18:   111 -> 114
#   F20: 11103 -> 11104
#   F21: 10864 -> 10865

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11526
    3  = 11008
    6  = 998
    7  = 117
    9  = 40
    10 = 38
    11 = 10640
    12 = 4113
    15 = 2457
    18 = 114
    20 = 11104
    21 = 10865
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
